$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 2.62

# Row 4
$ws.Range("G4").Value = 2.37
$ws.Range("H4").Value = 3.05
$ws.Range("J4").Value = 2.92
$ws.Range("L4").Value = 3.5
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 6.8
$ws.Range("O4").Value = 1.33
$ws.Range("P4").Value = 3.05
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.78
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.72
$ws.Range("U4").Value = 1.72
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 8.25
$ws.Range("X4").Value = 12.5
$ws.Range("Z4").Value = 26
$ws.Range("AA4").Value = 19
$ws.Range("AB4").Value = 27
$ws.Range("AC4").Value = 6.8
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 60
$ws.Range("AG4").Value = 8.5
$ws.Range("AH4").Value = 14.5
$ws.Range("AI4").Value = 10.5
$ws.Range("AK4").Value = 26
$ws.Range("AL4").Value = 35
$ws.Range("AM4").Value = 450
$ws.Range("AN4").Value = 4.4
$ws.Range("AO4").Value = 12.5
$ws.Range("AP4").Value = 18.5
$ws.Range("AQ4").Value = 50
$ws.Range("AR4").Value = 75
$ws.Range("AS4").Value = 200
$ws.Range("AT4").Value = 2.72
$ws.Range("AU4").Value = 6.7
$ws.Range("AV4").Value = 55
$ws.Range("AW4").Value = 4.9
$ws.Range("AX4").Value = 16
$ws.Range("AY4").Value = 23
$ws.Range("AZ4").Value = 75
$ws.Range("BA4").Value = 110
$ws.Range("BB4").Value = 300

# Row 5
$ws.Range("G5").Value = 8.75
$ws.Range("H5").Value = 4.55
$ws.Range("I5").Value = 1.3
$ws.Range("J5").Value = 7.8
$ws.Range("L5").Value = 1.78
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 3.85
$ws.Range("Q5").Value = 1.7
$ws.Range("R5").Value = 2.07
$ws.Range("S5").Value = 1.34
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 2.05
$ws.Range("W5").Value = 23
$ws.Range("X5").Value = 65
$ws.Range("Y5").Value = 28
$ws.Range("Z5").Value = 300
$ws.Range("AA5").Value = 120
$ws.Range("AB5").Value = 100
$ws.Range("AD5").Value = 9.5
$ws.Range("AG5").Value = 6.8
$ws.Range("AH5").Value = 6.1
$ws.Range("AJ5").Value = 8
$ws.Range("AK5").Value = 11
$ws.Range("AN5").Value = 9.5
$ws.Range("AO5").Value = 55
$ws.Range("AP5").Value = 50
$ws.Range("AQ5").Value = 450
$ws.Range("AR5").Value = 450
$ws.Range("AT5").Value = 3
$ws.Range("AW5").Value = 3.1
$ws.Range("AX5").Value = 5.9
$ws.Range("AZ5").Value = 16
$ws.Range("BA5").Value = 45

# Row 6
$ws.Range("G6").Value = 1.29
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 1.72
$ws.Range("K6").Value = 2.6
$ws.Range("L6").Value = 7
$ws.Range("N6").Value = 9.5
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 4.75
$ws.Range("Q6").Value = 1.47
$ws.Range("R6").Value = 2.5
$ws.Range("S6").Value = 1.27
$ws.Range("T6").Value = 3.4
$ws.Range("U6").Value = 1.82
$ws.Range("V6").Value = 1.9
$ws.Range("W6").Value = 9
$ws.Range("X6").Value = 7.2
$ws.Range("Y6").Value = 8.75
$ws.Range("Z6").Value = 8.5
$ws.Range("AA6").Value = 10.25
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 9.5
$ws.Range("AD6").Value = 10.75
$ws.Range("AE6").Value = 20
$ws.Range("AF6").Value = 75
$ws.Range("AG6").Value = 27
$ws.Range("AH6").Value = 65
$ws.Range("AI6").Value = 26
$ws.Range("AJ6").Value = 200
$ws.Range("AK6").Value = 90
$ws.Range("AL6").Value = 70
$ws.Range("AM6").Value = 500
$ws.Range("AN6").Value = 3.3
$ws.Range("AO6").Value = 5.7
$ws.Range("AQ6").Value = 14
$ws.Range("AR6").Value = 35
$ws.Range("AT6").Value = 3.4
$ws.Range("AU6").Value = 8.25
$ws.Range("AV6").Value = 70
$ws.Range("AW6").Value = 9.5
$ws.Range("AX6").Value = 45
$ws.Range("AY6").Value = 40
$ws.Range("AZ6").Value = 350
$ws.Range("BA6").Value = 300
$ws.Range("BB6").Value = 500

